# Adds two new columns, I ("I0") and J ("IF"), with per-row numeric data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, centered alignment) used by the
# other header cells (e.g. H1) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (2 .. 64) for columns I and J ---
$data = @(
    @(2,7,8), @(3,8,8), @(4,7,7), @(5,6,6), @(6,7,7),
    @(7,8,9), @(8,8,8), @(9,7,8), @(10,7,8), @(11,7,7),
    @(12,8,9), @(13,7,8), @(14,7,8), @(15,8,8), @(16,7,8),
    @(17,8,8), @(18,6,7), @(19,7,8), @(20,8,9), @(21,6,7),
    @(22,6,7), @(23,5,6), @(24,8,8), @(25,6,8), @(26,7,8),
    @(27,6,6), @(28,11,11), @(29,4,6), @(30,9,9), @(31,7,8),
    @(32,8,8), @(33,7,7), @(34,9,10), @(35,10,10), @(36,6,6),
    @(37,9,9), @(38,6,7), @(39,8,8), @(40,9,9), @(41,7,7),
    @(42,7,8), @(43,8,8), @(44,8,8), @(45,8,8), @(46,7,7),
    @(47,6,6), @(48,11,11), @(49,6,6), @(50,9,9), @(51,7,7),
    @(52,5,5), @(53,8,8), @(54,8,8), @(55,8,8), @(56,6,6),
    @(57,8,8), @(58,8,8), @(59,7,7), @(60,7,7), @(61,7,7),
    @(62,7,7), @(63,7,7), @(64,8,8)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
